$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style donor cell: a plain unstyled text cell (column B data rows), used to
# reset the style of cells that need an explicit text (quote-prefix) coercion
# back to the workbook default (style 0) after forcing a numeric-looking string
# to be stored as text.
$plainStyle = $ws.Cells.Item(2, 2).Style

# Row 2
$ws.Cells.Item(2, 4).Value = "40.218.17"
$ws.Cells.Item(2, 5).Value = "  +0.21%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.213.38"
$ws.Cells.Item(3, 5).Value = "  -0.57%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.01%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'295.88"
$ws.Cells.Item(5, 4).Style = $plainStyle
$ws.Cells.Item(5, 5).Value = "  +1.31%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'88.02"
$ws.Cells.Item(6, 4).Style = $plainStyle
$ws.Cells.Item(6, 5).Value = "  +0.35%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "'0.514"
$ws.Cells.Item(7, 4).Style = $plainStyle
$ws.Cells.Item(7, 5).Value = "  +0.41%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.07%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  -0.47%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "'52.17"
$ws.Cells.Item(10, 4).Style = $plainStyle
$ws.Cells.Item(10, 5).Value = "  +7.29%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "'30.97"
$ws.Cells.Item(11, 4).Style = $plainStyle
$ws.Cells.Item(11, 5).Value = "  +1.85%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "'0.0782"
$ws.Cells.Item(12, 4).Style = $plainStyle
$ws.Cells.Item(12, 5).Value = "  -0.07%  "

# Row 13
$ws.Cells.Item(13, 5).Value = "  +2.47%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "'6.40"
$ws.Cells.Item(14, 4).Style = $plainStyle
$ws.Cells.Item(14, 5).Value = "  -1.07%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "2.554.69"
$ws.Cells.Item(15, 5).Value = "  -0.67%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "'13.85"
$ws.Cells.Item(16, 4).Style = $plainStyle
$ws.Cells.Item(16, 5).Value = "  -0.81%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "2.213.36"
$ws.Cells.Item(17, 5).Value = "  -0.33%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "'0.737"
$ws.Cells.Item(18, 4).Style = $plainStyle
$ws.Cells.Item(18, 5).Value = "  +0.93%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "40.125.62"
$ws.Cells.Item(19, 5).Value = "  +0.12%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "0.0₃0888"
$ws.Cells.Item(20, 5).Value = "  -0.09%  "

# Row 21
$ws.Cells.Item(21, 5).Value = "  -0.43%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'5.77"
$ws.Cells.Item(22, 4).Style = $plainStyle
$ws.Cells.Item(22, 5).Value = "  -1.00%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "'65.71"
$ws.Cells.Item(23, 4).Style = $plainStyle
$ws.Cells.Item(23, 5).Value = "  +0.05%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "'235.86"
$ws.Cells.Item(24, 4).Style = $plainStyle
$ws.Cells.Item(24, 5).Value = "  -0.39%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  -0.03%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'2.49"
$ws.Cells.Item(26, 4).Style = $plainStyle
$ws.Cells.Item(26, 5).Value = "  +1.02%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  -0.84%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "'23.24"
$ws.Cells.Item(28, 4).Style = $plainStyle
$ws.Cells.Item(28, 5).Value = "  +2.21%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "'9.33"
$ws.Cells.Item(29, 4).Style = $plainStyle
$ws.Cells.Item(29, 5).Value = "  +1.05%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  -4.95%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "'156.57"
$ws.Cells.Item(31, 4).Style = $plainStyle
$ws.Cells.Item(31, 5).Value = "  +0.18%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "'32.20"
$ws.Cells.Item(32, 4).Style = $plainStyle
$ws.Cells.Item(32, 5).Value = "  +0.98%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  +0.04%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "'4.98"
$ws.Cells.Item(34, 4).Style = $plainStyle
$ws.Cells.Item(34, 5).Value = "  +0.35%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "'3.01"
$ws.Cells.Item(35, 4).Style = $plainStyle
$ws.Cells.Item(35, 5).Value = "  +3.42%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "'0.0716"
$ws.Cells.Item(36, 4).Style = $plainStyle
$ws.Cells.Item(36, 5).Value = "  -0.60%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "'0.114"
$ws.Cells.Item(38, 4).Style = $plainStyle
$ws.Cells.Item(38, 5).Value = "  +1.75%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  +2.89%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  +1.98%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "'15.62"
$ws.Cells.Item(41, 4).Style = $plainStyle
$ws.Cells.Item(41, 5).Value = "  -1.03%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "'3.82"
$ws.Cells.Item(42, 4).Style = $plainStyle
$ws.Cells.Item(42, 5).Value = "  -1.35%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "2.071.42"
$ws.Cells.Item(43, 5).Value = "  -2.61%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "'19.37"
$ws.Cells.Item(44, 4).Style = $plainStyle
$ws.Cells.Item(44, 5).Value = "  +5.72%  "

# Row 45
$ws.Cells.Item(45, 5).Value = "  +0.79%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "'9.99"
$ws.Cells.Item(46, 4).Style = $plainStyle
$ws.Cells.Item(46, 5).Value = "  -0.09%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "'2.82"
$ws.Cells.Item(47, 4).Style = $plainStyle
$ws.Cells.Item(47, 5).Value = "  +5.37%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "'1.90"
$ws.Cells.Item(48, 4).Style = $plainStyle
$ws.Cells.Item(48, 5).Value = "  -11.53%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "2.427.36"
$ws.Cells.Item(49, 5).Value = "  -0.39%  "

# Row 50
$ws.Cells.Item(50, 5).Value = "  +1.75%  "

# Row 51
$ws.Cells.Item(51, 5).Value = "  +0.76%  "

